# Re-upload of the workbook with refreshed calculation results.
# (mirrors the author's "Add files via upload" commit: the numeric
# results in the three data sheets were recomputed/replaced.)

$wb = $excel.ActiveWorkbook

# Helper: force a numeric-looking value to be stored as literal TEXT
# (t="s" shared-string cell), the way Excel does internally when a
# cell is formatted as Text before the value is typed in. A straight
# "$range.Value = '0.123'" assignment would otherwise get silently
# re-parsed back into a Number.
function Set-TextValue {
    param($range, [string]$text)

    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
}

# ---- Sheet 2: "Penambahan PDRB, LT, PDD" ------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B2").Value = "0.28438578293359296"
Set-TextValue $ws2.Range("C2") "0.4522788374012998"
$ws2.Range("D2").Value = -790448570037472

$ws2.Range("B7").Value = "0.316882237134571"
$ws2.Range("B8").Value = "0.4796883930650522"
$ws2.Range("B9").Value = "0.039839962742649115"

# ---- Sheet 3: "PDRB Perkotaan 2010, LT, PDD" ---------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B2").Value = "0.28438578293359296"
Set-TextValue $ws3.Range("C2") "0.4522788374012998"
$ws3.Range("D2").Value = -790448570037472

$ws3.Range("B6").Value = "0.316882237134571"

# These two values are dropped entirely in the refreshed workbook.
$ws3.Range("B7").ClearContents()
$ws3.Range("B8").ClearContents()

# Sheet 1 ("PeningkProduktivitasEfisiensi") and Sheet 4
# ("PDRB Perkotaan 2020, LT, PDD") keep the same data values in the
# refreshed file, so nothing to change there.
